$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$d2 = $ws.Range("D2")
$d2.NumberFormat = "@"
$d2.Value = "29.816.64"
$d2.Style = "Normal"
$ws.Range("E2").Value = "  -0.40%  "

$d3 = $ws.Range("D3")
$d3.NumberFormat = "@"
$d3.Value = "1.893.01"
$d3.Style = "Normal"
$ws.Range("E3").Value = "  +0.10%  "

$d4 = $ws.Range("D4")
$d4.NumberFormat = "@"
$d4.Value = "1.001"
$d4.Style = "Normal"
$ws.Range("E4").Value = "  -0.04%  "

$d5 = $ws.Range("D5")
$d5.NumberFormat = "@"
$d5.Value = "0.7952"
$d5.Style = "Normal"
$ws.Range("E5").Value = "  -3.20%  "

$d6 = $ws.Range("D6")
$d6.NumberFormat = "@"
$d6.Value = "242.90"
$d6.Style = "Normal"
$ws.Range("E6").Value = "  +0.59%  "

$ws.Range("E7").Value = "  +0.02%  "

$d8 = $ws.Range("D8")
$d8.NumberFormat = "@"
$d8.Value = "0.3165"
$d8.Style = "Normal"
$ws.Range("E8").Value = "  -2.51%  "

$ws.Range("E9").Value = "  -3.79%  "

$d10 = $ws.Range("D10")
$d10.NumberFormat = "@"
$d10.Value = "0.07060"
$d10.Style = "Normal"
$ws.Range("E10").Value = "  +0.49%  "

$d11 = $ws.Range("D11")
$d11.NumberFormat = "@"
$d11.Value = "0.08081"
$d11.Style = "Normal"
$ws.Range("E11").Value = "  +0.58%  "

$d12 = $ws.Range("D12")
$d12.NumberFormat = "@"
$d12.Value = "0.7677"
$d12.Style = "Normal"
$ws.Range("E12").Value = "  +2.95%  "

$d13 = $ws.Range("D13")
$d13.NumberFormat = "@"
$d13.Value = "1.906.02"
$d13.Style = "Normal"
$ws.Range("E13").Value = "  +0.79%  "

$d14 = $ws.Range("D14")
$d14.NumberFormat = "@"
$d14.Value = "5.354"
$d14.Style = "Normal"
$ws.Range("E14").Value = "  +2.89%  "

$d15 = $ws.Range("D15")
$d15.NumberFormat = "@"
$d15.Value = "92.51"
$d15.Style = "Normal"
$ws.Range("E15").Value = "  +0.37%  "

$d16 = $ws.Range("D16")
$d16.NumberFormat = "@"
$d16.Value = "29.835.42"
$d16.Style = "Normal"
$ws.Range("E16").Value = "  -0.29%  "

$d17 = $ws.Range("D17")
$d17.NumberFormat = "@"
$d17.Value = "6.006"
$d17.Style = "Normal"
$ws.Range("E17").Value = "  +2.04%  "

$ws.Range("E18").Value = "  -1.14%  "

$d19 = $ws.Range("D19")
$d19.NumberFormat = "@"
$d19.Value = "244.46"
$d19.Style = "Normal"
$ws.Range("E19").Value = "  -0.10%  "

$d20 = $ws.Range("D20")
$d20.NumberFormat = "@"
$d20.Value = "0.000007710"
$d20.Style = "Normal"
$ws.Range("E20").Value = "  -0.58%  "

$d21 = $ws.Range("D21")
$d21.NumberFormat = "@"
$d21.Value = "8.351"
$d21.Style = "Normal"
$ws.Range("E21").Value = "  +20.63%  "

$d23 = $ws.Range("D23")
$d23.NumberFormat = "@"
$d23.Value = "2.150.43"
$d23.Style = "Normal"
$ws.Range("E23").Value = "  +0.52%  "

$ws.Range("E24").Value = "  -0.03%  "

$d25 = $ws.Range("D25")
$d25.NumberFormat = "@"
$d25.Value = "0.1639"
$d25.Style = "Normal"
$ws.Range("E25").Value = "  +5.65%  "

$d26 = $ws.Range("D26")
$d26.NumberFormat = "@"
$d26.Value = "9.358"
$d26.Style = "Normal"
$ws.Range("E26").Value = "  +1.82%  "

$d27 = $ws.Range("D27")
$d27.NumberFormat = "@"
$d27.Value = "166.20"
$d27.Style = "Normal"
$ws.Range("E27").Value = "  +0.10%  "

$ws.Range("E28").Value = "  -0.47%  "

$ws.Range("E29").Value = "  -1.61%  "

$d30 = $ws.Range("D30")
$d30.NumberFormat = "@"
$d30.Value = "1.399"
$d30.Style = "Normal"
$ws.Range("E30").Value = "  +2.36%  "

$ws.Range("E31").Value = "  +1.40%  "

$d32 = $ws.Range("D32")
$d32.NumberFormat = "@"
$d32.Value = "4.437"
$d32.Style = "Normal"
$ws.Range("E32").Value = "  +3.95%  "

$d33 = $ws.Range("D33")
$d33.NumberFormat = "@"
$d33.Value = "0.05709"
$d33.Style = "Normal"
$ws.Range("E33").Value = "  +1.42%  "

$d34 = $ws.Range("D34")
$d34.NumberFormat = "@"
$d34.Value = "4.046"
$d34.Style = "Normal"
$ws.Range("E34").Value = "  -0.52%  "

$d35 = $ws.Range("D35")
$d35.NumberFormat = "@"
$d35.Value = "1.261"
$d35.Style = "Normal"
$ws.Range("E35").Value = "  -0.68%  "

$d36 = $ws.Range("D36")
$d36.NumberFormat = "@"
$d36.Value = "0.7390"
$d36.Style = "Normal"
$ws.Range("E36").Value = "  +1.41%  "

$d37 = $ws.Range("D37")
$d37.NumberFormat = "@"
$d37.Value = "0.9986"
$d37.Style = "Normal"
$ws.Range("E37").Value = "  -0.10%  "

$d38 = $ws.Range("D38")
$d38.NumberFormat = "@"
$d38.Value = "2.630"
$d38.Style = "Normal"
$ws.Range("E38").Value = "  -3.03%  "

$d39 = $ws.Range("D39")
$d39.NumberFormat = "@"
$d39.Value = "0.01910"
$d39.Style = "Normal"
$ws.Range("E39").Value = "  -0.09%  "

$d40 = $ws.Range("D40")
$d40.NumberFormat = "@"
$d40.Value = "2.784"
$d40.Style = "Normal"
$ws.Range("E40").Value = "  +0.14%  "

$d41 = $ws.Range("D41")
$d41.NumberFormat = "@"
$d41.Value = "0.4407"
$d41.Style = "Normal"
$ws.Range("E41").Value = "  -0.32%  "

$d42 = $ws.Range("D42")
$d42.NumberFormat = "@"
$d42.Value = "72.53"
$d42.Style = "Normal"
$ws.Range("E42").Value = "  +0.94%  "

$d43 = $ws.Range("D43")
$d43.NumberFormat = "@"
$d43.Value = "5.814"
$d43.Style = "Normal"
$ws.Range("E43").Value = "  -2.45%  "

$d44 = $ws.Range("D44")
$d44.NumberFormat = "@"
$d44.Value = "0.8413"
$d44.Style = "Normal"
$ws.Range("E44").Value = "  -0.16%  "

$ws.Range("E45").Value = "  +0.05%  "

$d46 = $ws.Range("D46")
$d46.NumberFormat = "@"
$d46.Value = "1.034.53"
$d46.Style = "Normal"
$ws.Range("E46").Value = "  +4.41%  "

$d47 = $ws.Range("D47")
$d47.NumberFormat = "@"
$d47.Value = "103.14"
$d47.Style = "Normal"
$ws.Range("E47").Value = "  +2.58%  "

$d48 = $ws.Range("D48")
$d48.NumberFormat = "@"
$d48.Value = "1.874"
$d48.Style = "Normal"
$ws.Range("E48").Value = "  +0.07%  "

$d49 = $ws.Range("D49")
$d49.NumberFormat = "@"
$d49.Value = "9.984"
$d49.Style = "Normal"
$ws.Range("E49").Value = "  +2.47%  "

$d50 = $ws.Range("D50")
$d50.NumberFormat = "@"
$d50.Value = "7.428"
$d50.Style = "Normal"
$ws.Range("E50").Value = "  -1.89%  "

$d51 = $ws.Range("D51")
$d51.NumberFormat = "@"
$d51.Value = "2.041.09"
$d51.Style = "Normal"
$ws.Range("E51").Value = "  +0.10%  "
